$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full target product data for rows 2..17 (A..F)
$data = @(
    @("Nothing Ear (2)", 135, 146.15, 8.259259259259263, 11.15000000000001, "https://www.idealo.de/preisvergleich/OffersOfProduct/202568190_-ear-2-nothing-tech.html"),
    @("Huawei Freebuds 5i", 71.09999999999999, 78, 9.704641350210979, 6.900000000000006, "https://www.idealo.de/preisvergleich/OffersOfProduct/202282533_-freebuds-5i-huawei.html"),
    @("LG TONE Free DFP8", 59.4, 69, 16.16161616161617, 9.600000000000001, "https://www.idealo.de/preisvergleich/OffersOfProduct/201653276_-tone-free-dfp8-lg-electronics.html"),
    @("Apple iPhone 6S", 119.95, 248.08, 106.8195081283868, 128.13, "https://www.idealo.de/preisvergleich/OffersOfProduct/4842805_-iphone-6s-apple.html"),
    @("Bowers & Wilkins PX7 S2", 199, 428.37, 115.2613065326633, 229.37, "https://www.idealo.de/preisvergleich/OffersOfProduct/202011022_-px7-s2-schwarz-bowers-wilkins.html"),
    @("Sony WH-1000XM5", 296.57, 299, 0.8193681087095818, 2.430000000000007, "https://www.idealo.de/preisvergleich/OffersOfProduct/201952086_-wh-1000xm5-black-sony.html"),
    @("Sony WH-1000XM4", 249.05, 269.99, 8.407950210801042, 20.94, "https://www.idealo.de/preisvergleich/OffersOfProduct/200278050_-wh-1000xm4-schwarz-sony.html"),
    @("JBL Wave 200TWS", 42.77, 52.48, 22.70282908580779, 9.709999999999994, "https://www.idealo.de/preisvergleich/OffersOfProduct/201746578_-wave-200tws-jbl-audio.html"),
    @("Sony WH-CH720N", 79.98999999999999, 86.98, 8.738592324040518, 6.990000000000009, "https://www.idealo.de/preisvergleich/OffersOfProduct/202334761_-wh-ch720n-sony.html"),
    @("Sony LinkBuds S", 129, 146, 13.17829457364341, 17, "https://www.idealo.de/preisvergleich/OffersOfProduct/201961566_-linkbuds-s-schwarz-sony.html"),
    @("Marshall Major IV", 79.98999999999999, 129.9, 62.39529941242657, 49.91000000000001, "https://www.idealo.de/preisvergleich/OffersOfProduct/200799717_-major-iv-marshall.html"),
    @("LG TONE Free DFP9E", 76.48999999999999, 84.98999999999999, 11.11256373382142, 8.5, "https://www.idealo.de/preisvergleich/OffersOfProduct/202056848_-tone-free-dfp9e-lg-electronics.html"),
    @("Samsung Galaxy Buds 2 Pro SM-R510", 139.9, 169.09, 20.86490350250178, 29.19, "https://www.idealo.de/preisvergleich/OffersOfProduct/202054812_-galaxy-buds-2-pro-sm-r510-white-samsung.html"),
    @("Skullcandy Crusher Evo", 122.8, 144.99, 18.07003257328991, 22.19000000000001, "https://www.idealo.de/preisvergleich/OffersOfProduct/200818290_-crusher-evo-skullcandy.html"),
    @("Beats By Dre Studio3 Wireless", 197.91, 269, 35.92036784396948, 71.09, "https://www.idealo.de/preisvergleich/OffersOfProduct/5760160_-studio3-wireless-beats-by-dr-dre.html"),
    @("Xiaomi Buds 3T Pro", 35.99, 54.43, 51.23645457071407, 18.44, "https://www.idealo.de/preisvergleich/OffersOfProduct/201838960_-buds-3t-pro-xiaomi.html")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $item = $data[$i]
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    $ws.Cells.Item($row, 4).Value = $item[3]
    $ws.Cells.Item($row, 5).Value = $item[4]
    $ws.Cells.Item($row, 6).Value = $item[5]
}
